# Update "想去人数" (attendance count) figures across all sheets to match
# the refreshed data snapshot (gh-pages output regenerated at 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 26
$ws.Range("F7").Value = 1136
$ws.Range("F8").Value = 366
$ws.Range("F9").Value = 226
$ws.Range("F10").Value = 321
$ws.Range("F11").Value = 7996
$ws.Range("F13").Value = 9477
$ws.Range("F14").Value = 74
$ws.Range("F16").Value = 8
$ws.Range("F17").Value = 464
$ws.Range("F18").Value = 5
$ws.Range("F27").Value = 380
$ws.Range("F29").Value = 1622
$ws.Range("F30").Value = 25
$ws.Range("F31").Value = 65
$ws.Range("F32").Value = 300
$ws.Range("F35").Value = 336
$ws.Range("F36").Value = 56
$ws.Range("F37").Value = 929
$ws.Range("F42").Value = 306
$ws.Range("F45").Value = 272
$ws.Range("F47").Value = 239
$ws.Range("F48").Value = 89

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 93
$ws.Range("G9").Value = 108
$ws.Range("F13").Value = 12
$ws.Range("F19").Value = 18

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 186

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 186
$ws.Range("F11").Value = 1136
$ws.Range("F12").Value = 366
$ws.Range("F15").Value = 93
$ws.Range("F16").Value = 321
$ws.Range("F17").Value = 7996
$ws.Range("F18").Value = 9477
$ws.Range("F19").Value = 75
$ws.Range("F20").Value = 8
$ws.Range("F21").Value = 464
$ws.Range("F27").Value = 380
$ws.Range("F28").Value = 1622
$ws.Range("F29").Value = 25
$ws.Range("F30").Value = 65
$ws.Range("F31").Value = 300
$ws.Range("F33").Value = 336
$ws.Range("F34").Value = 56
$ws.Range("G35").Value = 108
$ws.Range("F36").Value = 929
$ws.Range("F38").Value = 12
$ws.Range("F41").Value = 306
$ws.Range("F44").Value = 272
$ws.Range("F46").Value = 239
$ws.Range("F47").Value = 18
$ws.Range("F49").Value = 89
